$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new data row (row 8) ---
# Copy formatting (date number format) from the cell above (A7) so the new
# date cell picks up the same style index used by the rest of column A.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A8").Value = 46054
$ws.Range("B8").Value = "Cadel Evans Great Ocean Road Race"
$ws.Range("C8").Value = "One day race"
$ws.Range("D8").Value = "Tobias Lund Andresen"
$ws.Range("E8").Value = "Matthew Brennan"
$ws.Range("F8").Value = "Brady Gilmore"
$ws.Range("G8").Value = "Mauro Schmid"
$ws.Range("H8").Value = "Natnael Tesfatsion"
$ws.Range("I8").Value = "Laurence Pithie"
$ws.Range("J8").Value = "Filippo Zana"
$ws.Range("K8").Value = "Gal Glivar"
$ws.Range("L8").Value = "Francesco Busatto"
$ws.Range("M8").Value = "Aaron Gate"

# --- Resize columns to their (new) best-fit widths now that the longer
#     strings above have been added ---
$ws.Columns.Item(1).ColumnWidth = 8.833333333333334
$ws.Columns.Item(2).ColumnWidth = 29.833333333333332
$ws.Columns.Item(3).ColumnWidth = 10.666666666666666
$ws.Columns.Item(4).ColumnWidth = 18.333333333333332
$ws.Columns.Item(5).ColumnWidth = 18.333333333333332
$ws.Columns.Item(6).ColumnWidth = 18.333333333333332
$ws.Columns.Item(7).ColumnWidth = 14.666666666666666
$ws.Columns.Item(8).ColumnWidth = 15.333333333333334
$ws.Columns.Item(9).ColumnWidth = 16.0
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 16.5
$ws.Columns.Item(12).ColumnWidth = 20.666666666666668
$ws.Columns.Item(13).ColumnWidth = 13.833333333333334

# --- Update the view: zoom to 94%, scroll back to A1, and select the
#     whole sheet (mirrors an "A1:XFD1048576" selection) ---
$win = $excel.ActiveWindow
$win.Zoom = 94
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Cells.Select()
